# Update the cryptos list (Price / Volume(1h) columns) with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds values that look numeric (e.g. "243.00", "0.9983") but
# must be preserved verbatim as text (trailing zeros, scientific-looking small
# decimals, etc.). Force the column to Text format first so Excel does not
# silently convert/round these values when we assign them.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 43 and 44 swap coin identity (Aave moves above FraxShare) plus new D/E values.
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "72.81"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "6.038"
$ws.Range("E44").Value = "  -0.74%  "

# Map of row -> (Price, Volume(1h)) values for all other rows that changed.
$updates = @{
    2  = @("30.036.73", "  +0.62%  ")
    3  = @("1.889.41", "  +0.18%  ")
    4  = @("0.9983", "  -0.28%  ")
    5  = @("0.7455", "  -2.61%  ")
    6  = @("243.00", "  +0.26%  ")
    7  = @("0.9969", "  -0.39%  ")
    8  = @("0.3169", "  +1.40%  ")
    9  = @("0.07252", "  +1.68%  ")
    10 = @("25.06", "  -1.82%  ")
    11 = @("0.08360", "  -1.87%  ")
    12 = @("0.7600", "  -0.14%  ")
    13 = @("5.426", "  +1.34%  ")
    14 = @("1.819.80", "  -4.29%  ")
    15 = @("92.68", "  -0.81%  ")
    16 = @("6.159", "  +0.61%  ")
    17 = @("29.982.17", "  +0.19%  ")
    18 = @("249.93", "  +2.58%  ")
    19 = @("13.65", "  -0.45%  ")
    20 = @("0.000007877", "  +0.85%  ")
    21 = @("2.160.59", "  +0.68%  ")
    23 = @("8.002", "  -0.08%  ")
    24 = @("0.9982", "  -0.30%  ")
    25 = @("0.1585", "  -2.11%  ")
    26 = @("9.312", "  -0.70%  ")
    27 = @("165.43", "  +1.12%  ")
    28 = @("18.77", "  +0.05%  ")
    29 = @("2.050", "  +0.92%  ")
    30 = @("1.489", "  -1.38%  ")
    31 = @("4.616", "  +2.68%  ")
    32 = @("1.538", "  -0.05%  ")
    33 = @("4.234", "  +2.98%  ")
    34 = @("0.05384", "  -0.95%  ")
    35 = @("1.255", "  +1.26%  ")
    36 = @("0.7592", "  +2.04%  ")
    37 = @("0.9906", "  -1.02%  ")
    38 = @("2.709", "  +0.48%  ")
    39 = @("0.01969", "  +1.20%  ")
    40 = @("2.767", "  -0.43%  ")
    41 = @("0.4579", "  +2.74%  ")
    42 = @("1.104.21", "  +0.35%  ")
    45 = @("0.8716", "  +1.55%  ")
    46 = @("104.54", "  +1.91%  ")
    47 = @("0.9992", "  -0.13%  ")
    48 = @("1.872", "  +0.34%  ")
    49 = @("7.630", "  -0.25%  ")
    50 = @("9.605", "  -1.26%  ")
    51 = @("2.034.40", "  -0.44%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

# Row 22 only has its Price updated; Volume(1h) stays the same.
$ws.Range("D22").Value = "0.9981"
